$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.816.19'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.635.43'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.37'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0642'
$ws.Range('E10').Value = '  +1.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0781'
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.640.97'
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.860.54'
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.557'
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₃0770'
$ws.Range('E16').Value = '  +1.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.07'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.827.01'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '194.12'
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('E22').Value = '  +0.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.14'
$ws.Range('E23').Value = '  +1.57%  '
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('E25').Value = '  -1.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '139.63'
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('E27').Value = '  -5.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.83'
$ws.Range('E28').Value = '  +1.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.53'
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0497'
$ws.Range('E31').Value = '  +1.86%  '
$ws.Range('E32').Value = '  +1.26%  '
$ws.Range('E33').Value = '  +1.21%  '
$ws.Range('E34').Value = '  +2.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.39'
$ws.Range('E35').Value = '  +0.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.903'
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.554'
$ws.Range('E38').Value = '  +0.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.111.70'
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('E42').Value = '  +0.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.42'
$ws.Range('E43').Value = '  +2.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.800'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0₆0112'
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('E46').Value = '  +13.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.51'
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.419'
$ws.Range('E48').Value = '  -5.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.68'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('E50').Value = '  -0.34%  '
$ws.Range('E51').Value = '  +0.52%  '
